$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("B1").Value = "Weaty FC"
$ws.Range("B2").Value = "Oakridge Nomads"
$ws.Range("B3").Value = "Sunday 16 Nov"
$ws.Range("B5").Value = "Brighton Hill Playing Fields 2"
$ws.Range("B6").Value = "Bernard Cornish Cup"
